# commit: add Gui of user storie 30 add func 30 to the main game
$wb = $excel.ActiveWorkbook

# --- Cards sheet: flip the "enable" flags (user story 30) ---
$cards = $wb.Worksheets.Item("Cards")
$cards.Range("I3").Value = "no"
$cards.Range("I9").Value = "yes"
$cards.Range("I10").Value = "yes"

# --- Games sheet: record the new feedback for game/row 11, drop the old bulk rows ---
$games = $wb.Worksheets.Item("Games")

$games.Range("B12").Value = 2
$games.Range("C12").Value = 6
$games.Range("D12").Value = 1

# Force the date to be stored as literal text (matches existing "date" column
# cells, which are plain text rather than real Excel dates).
$games.Range("E12").NumberFormat = "@"
$games.Range("E12").Value = "12.01.2020"
$games.Range("E12").ClearFormats()

$games.Range("F12").Value = "was ok"

# Rows 13-44 no longer apply after the new GUI/feedback flow - remove them so
# the sheet ends at row 12.
$games.Range("A13:A44").EntireRow.Delete()

# Leave the Games tab focused near the new row, like the author's view.
$games.Activate()
$null = $games.Range("I11").Select()
